# Split the "Describe Kafka ..." bullet on the exam-preview questions slide
# into three runs, dropping "logical, " from the sentence:
#   "Describe Kafka at a conceptual, logical, system, and physical level"
# becomes the concatenation of three runs:
#   "Describe Kafka at " + "a conceptual, " + "system, and physical level"

$p = $ppt.ActivePresentation

$targetOld  = "Describe Kafka at a conceptual, logical, system, and physical level"
$run1Text   = "Describe Kafka at "
$run2Text   = "a conceptual, "
$run3Text   = "system, and physical level"

$targetParagraph = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $paraCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $paraCount; $pi++) {
                    $candidate = $tr.Paragraphs($pi)
                    if ($candidate.Text -eq $targetOld) {
                        $targetParagraph = $candidate
                    }
                }
            }
        }
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not find the 'Describe Kafka ...' paragraph"
}

# Shrink the existing (only) run down to the first chunk, then append the
# other two chunks as their own runs via InsertAfter -- this reproduces the
# run split shown in the diff.
$firstRun = $targetParagraph.Runs(1)
$firstRun.Text = $run1Text
$targetParagraph.InsertAfter($run2Text) | Out-Null
$targetParagraph.InsertAfter($run3Text) | Out-Null
